# Populate the "D" (Test Result) columns for the test cases that were
# actually executed, for the three sheets touched by this commit:
#   - Item Checkout Test Cases : "Child Checkout Validation" block (rows 24-26)
#   - Return Item Test Cases   : "Return Item Validation" block   (rows 4-10)
#   - OverDueTestCases         : "Overdue item fee > item value"  block (rows 22-24)
#                                 "Paying overdue fines"           block (rows 41-47)
#
# The "Test Case Summary" sheet and each table's header row (e.g. D23, D3,
# D21, D40) are driven by formulas that reference these cells, so they
# recalculate automatically once the underlying results are entered.

$wb = $excel.ActiveWorkbook

# --- Item Checkout Test Cases: Child Checkout Validation ---
$ws = $wb.Worksheets.Item("Item Checkout Test Cases")
$ws.Range("D24").Value = "PASS"
$ws.Range("D25").Value = "PASS"
$ws.Range("D26").Value = "FAIL"

# --- Return Item Test Cases: Return Item Validation ---
$ws = $wb.Worksheets.Item("Return Item Test Cases")
$ws.Range("D4").Value = "PASS"
$ws.Range("D5").Value = "PASS"
$ws.Range("D6").Value = "PASS"
$ws.Range("D7").Value = "PASS"
$ws.Range("D8").Value = "PASS"
$ws.Range("D9").Value = "PASS"
$ws.Range("D10").Value = "PASS"

# --- OverDueTestCases: Overdue item fee > item value ---
$ws = $wb.Worksheets.Item("OverDueTestCases")
$ws.Range("D22").Value = "PASS"
$ws.Range("D23").Value = "PASS"
$ws.Range("D24").Value = "FAIL"

# --- OverDueTestCases: Paying overdue fines ---
$ws.Range("D41").Value = "PASS"
$ws.Range("D42").Value = "PASS"
$ws.Range("D43").Value = "PASS"
$ws.Range("D44").Value = "PASS"
$ws.Range("D45").Value = "PASS"
$ws.Range("D46").Value = "PASS"
$ws.Range("D47").Value = "PASS"
